$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 480-481),
# pushing the existing rows 480-498 down to 482-500.
$ws.Range("A480:R481").EntireRow.Insert()

# Row 480: new weekly record - Alcachofa, Argentina(o), Primera
$ws.Cells.Item(480, 1).Value = 3
$ws.Cells.Item(480, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(480, 3).Value = "Coquimbo"
$ws.Cells.Item(480, 4).Value = 45075
$ws.Cells.Item(480, 5).Value = 5
$ws.Cells.Item(480, 6).Value = 100112013
$ws.Cells.Item(480, 7).Value = "Alcachofa"
$ws.Cells.Item(480, 8).Value = "Argentina(o)"
$ws.Cells.Item(480, 9).Value = "Primera"
$ws.Cells.Item(480, 10).Value = 105
$ws.Cells.Item(480, 11).Value = 13500
$ws.Cells.Item(480, 12).Value = 14000
$ws.Cells.Item(480, 13).Value = 13738
$ws.Cells.Item(480, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(480, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(480, 16).Value = 275
$ws.Cells.Item(480, 17).Value = 50
$ws.Cells.Item(480, 18).Value = "Hortaliza"

# Row 481: new weekly record - Alcachofa, Española, Primera
$ws.Cells.Item(481, 1).Value = 3
$ws.Cells.Item(481, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(481, 3).Value = "Coquimbo"
$ws.Cells.Item(481, 4).Value = 45075
$ws.Cells.Item(481, 5).Value = 5
$ws.Cells.Item(481, 6).Value = 100112013
$ws.Cells.Item(481, 7).Value = "Alcachofa"
$ws.Cells.Item(481, 8).Value = "Española"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 165
$ws.Cells.Item(481, 11).Value = 16500
$ws.Cells.Item(481, 12).Value = 17000
$ws.Cells.Item(481, 13).Value = 16742
$ws.Cells.Item(481, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(481, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(481, 16).Value = 558
$ws.Cells.Item(481, 17).Value = 30
$ws.Cells.Item(481, 18).Value = "Hortaliza"
